$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New table (copy of the strategy/cost/qaly/ICER table, pasted as values
#     two columns to the right, starting at row 32) ---

# Header row
$ws.Range("G32").Value = "strategy"
$ws.Range("H32").Value = "final.cost"
$ws.Range("I32").Value = "final.qaly"
$ws.Range("J32").Value = "ICER"

# mc
$ws.Range("G33").Value = "mc"
$ws.Range("H33").Value = 24177.5805553879
$ws.Range("I33").Value = 0.61573595541761295
$ws.Range("J33").Value = "N/A"

# dc
$ws.Range("G34").Value = "dc"
$ws.Range("H34").Value = 24178.360498779901
$ws.Range("I34").Value = 0.61751482912255096
$ws.Range("J34").Value = 438.44787285111875

# edc
$ws.Range("G35").Value = "edc"
$ws.Range("H35").Value = 24241.862419323399
$ws.Range("I35").Value = 0.62051449850734697
$ws.Range("J35").Value = 21169.639849431824

# universal
$ws.Range("G36").Value = "universal"
$ws.Range("H36").Value = 24568.2377349134
$ws.Range("I36").Value = 0.62492680289312497
$ws.Range("J36").Value = 73969.356384839077

# Extra QALY-delta x1000 calculation below the table
$ws.Range("H38").Formula = "=I36-I33"
$ws.Range("H39").Formula = "=H38*1000"

# --- View state: scroll / zoom / selection ---
$ws.Activate()
[void]$ws.Range("H17").Select()
$excel.ActiveWindow.Zoom = 115
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
